# Auto-generated script to apply scheduled market-price data refresh
# to the Leve profit tables across all job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H80").Value = 1447.3636
$ws.Range("I80").Value = 1012
$ws.Range("J80").Value = 1969.8
$ws.Range("K80").Value = 3036
$ws.Range("L80").Value = 5909.4
$ws.Range("M80").Value = -2038
$ws.Range("N80").Value = -7905.4
$ws.Range("H83").Value = 1447.3636
$ws.Range("I83").Value = 1012
$ws.Range("J83").Value = 1969.8
$ws.Range("K83").Value = 9108
$ws.Range("L83").Value = 17728.2
$ws.Range("M83").Value = -4116
$ws.Range("N83").Value = -27712.2
$ws.Range("H86").Value = 174999.75
$ws.Range("I86").Value = 174999
$ws.Range("K86").Value = 174999
$ws.Range("M86").Value = -173876
$ws.Range("H88").Value = 14774
$ws.Range("J88").Value = 20447.6
$ws.Range("L88").Value = 20447.6
$ws.Range("N88").Value = -21259.6
$ws.Range("H89").Value = 174999.75
$ws.Range("I89").Value = 174999
$ws.Range("K89").Value = 874995
$ws.Range("M89").Value = -869379
$ws.Range("H91").Value = 14774
$ws.Range("J91").Value = 20447.6
$ws.Range("L91").Value = 20447.6
$ws.Range("N91").Value = -23255.6
$ws.Range("H100").Value = 5500
$ws.Range("I100").Value = 5500
$ws.Range("K100").Value = 5500
$ws.Range("M100").Value = -4959
$ws.Range("H135").Value = 828.6667
$ws.Range("I135").Value = 807.25
$ws.Range("K135").Value = 7265.25
$ws.Range("M135").Value = -4730.25
$ws.Range("H138").Value = 2635.9092
$ws.Range("I138").Value = 1998.75
$ws.Range("K138").Value = 5996.25
$ws.Range("M138").Value = -856.25
$ws.Range("H141").Value = 2348
$ws.Range("I141").Value = 2348
$ws.Range("K141").Value = 7044
$ws.Range("M141").Value = -1864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3613.7317
$ws.Range("I32").Value = 1267.5
$ws.Range("K32").Value = 1267.5
$ws.Range("M32").Value = -980.5
$ws.Range("H61").Value = 3332.6667
$ws.Range("I61").Value = 3332.6667
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3332.6667
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3120.6667
$ws.Range("N61").ClearContents()
$ws.Range("H132").Value = 1962
$ws.Range("I132").Value = 1962
$ws.Range("K132").Value = 5886
$ws.Range("M132").Value = -3356
$ws.Range("H136").Value = 3332.6667
$ws.Range("I136").Value = 3332.6667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9998.000100000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7448.000100000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6634.3335
$ws.Range("I86").Value = 5835.3335
$ws.Range("K86").Value = 5835.3335
$ws.Range("M86").Value = -4712.3335
$ws.Range("H89").Value = 6634.3335
$ws.Range("I89").Value = 5835.3335
$ws.Range("K89").Value = 29176.6675
$ws.Range("M89").Value = -23560.6675
$ws.Range("H99").Value = 1159.7778
$ws.Range("I99").Value = 1192.25
$ws.Range("K99").Value = 1192.25
$ws.Range("M99").Value = 305.75
$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 507.63635
$ws.Range("I16").Value = 439.42856
$ws.Range("K16").Value = 439.42856
$ws.Range("M16").Value = -152.42856
$ws.Range("H31").Value = 2108.8518
$ws.Range("I31").Value = 1647
$ws.Range("J31").Value = 2686.1667
$ws.Range("K31").Value = 1647
$ws.Range("L31").Value = 2686.1667
$ws.Range("M31").Value = -1352
$ws.Range("N31").Value = -3276.1667
$ws.Range("H34").Value = 2108.8518
$ws.Range("I34").Value = 1647
$ws.Range("J34").Value = 2686.1667
$ws.Range("K34").Value = 1647
$ws.Range("L34").Value = 2686.1667
$ws.Range("M34").Value = -1445
$ws.Range("N34").Value = -3090.1667
$ws.Range("H62").Value = 1500
$ws.Range("I62").Value = 1500
$ws.Range("K62").Value = 1500
$ws.Range("M62").Value = -876
$ws.Range("H65").Value = 1500
$ws.Range("I65").Value = 1500
$ws.Range("K65").Value = 7500
$ws.Range("M65").Value = -4380
$ws.Range("H113").Value = 507.63635
$ws.Range("I113").Value = 439.42856
$ws.Range("K113").Value = 439.42856
$ws.Range("M113").Value = 1730.57144
$ws.Range("I132").Value = 5500
$ws.Range("K132").Value = 16500
$ws.Range("M132").Value = -13970
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1309.6666
$ws.Range("I5").Value = 965
$ws.Range("J5").Value = 1999
$ws.Range("K5").Value = 2895
$ws.Range("L5").Value = 5997
$ws.Range("M5").Value = -2783
$ws.Range("N5").Value = -6221
$ws.Range("H117").Value = 503.375
$ws.Range("J117").Value = 623.8333
$ws.Range("L117").Value = 1871.4999
$ws.Range("N117").Value = -8755.499900000001
$ws.Range("H121").Value = 866.5
$ws.Range("J121").Value = 866.5
$ws.Range("L121").Value = 2599.5
$ws.Range("N121").Value = -5219.5
$ws.Range("H134").Value = 201118.8
$ws.Range("I134").Value = 201118.8
$ws.Range("K134").Value = 603356.3999999999
$ws.Range("M134").Value = -598286.3999999999
$ws.Range("H135").Value = 1309.6666
$ws.Range("I135").Value = 965
$ws.Range("J135").Value = 1999
$ws.Range("K135").Value = 8685
$ws.Range("L135").Value = 17991
$ws.Range("M135").Value = -6150
$ws.Range("N135").Value = -23061
$ws.Range("H137").Value = 1846.75
$ws.Range("I137").Value = 1629
$ws.Range("K137").Value = 4887
$ws.Range("M137").Value = 213
$ws.Range("H139").Value = 1835.1
$ws.Range("I139").Value = 1372.6666
$ws.Range("K139").Value = 4117.9998
$ws.Range("M139").Value = 1022.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 997.5
$ws.Range("I107").Value = 995
$ws.Range("K107").Value = 995
$ws.Range("M107").Value = 925
$ws.Range("H122").Value = 12806.4
$ws.Range("I122").Value = 16440.428
$ws.Range("J122").Value = 4327
$ws.Range("K122").Value = 49321.284
$ws.Range("L122").Value = 12981
$ws.Range("M122").Value = -46871.284
$ws.Range("N122").Value = -17881

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 88595
$ws.Range("I5").Value = 97750
$ws.Range("K5").Value = 97750
$ws.Range("M5").Value = -97638
$ws.Range("H107").Value = 1383.2727
$ws.Range("I107").Value = 912.8889
$ws.Range("K107").Value = 2738.6667
$ws.Range("M107").Value = -818.6667000000002
$ws.Range("H122").Value = 1483.5
$ws.Range("I122").Value = 1180.2
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3540.6
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1090.6
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 28375.36
$ws.Range("I126").Value = 29336.21
$ws.Range("J126").Value = 25332.666
$ws.Range("K126").Value = 88008.63
$ws.Range("L126").Value = 75997.99800000001
$ws.Range("M126").Value = -85538.63
$ws.Range("N126").Value = -80937.99800000001
